$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 448 (existing data shifts down by 2 rows,
# from 448:485 to 450:487).
$ws.Rows.Item(448).EntireRow.Insert()
$ws.Rows.Item(448).EntireRow.Insert()

# Fill in the new row 448 with the new weekly record.
$ws.Range("A448").Value = 11
$ws.Range("B448").Value = "Vega Monumental Concepción"
$ws.Range("C448").Value = "Bíobío"
$ws.Range("D448").Value = 45223
$ws.Range("E448").Value = 8
$ws.Range("F448").Value = 100114013
$ws.Range("G448").Value = "Zanahoria"
$ws.Range("H448").Value = "Sin especificar"
$ws.Range("I448").Value = "Primera"
$ws.Range("J448").Value = 250
$ws.Range("K448").Value = 5500
$ws.Range("L448").Value = 5500
$ws.Range("M448").Value = 5500
$ws.Range("N448").Value = "$/saco 20 kilos"
$ws.Range("O448").Value = "Región de Coquimbo"
$ws.Range("P448").Value = 275
$ws.Range("Q448").Value = 20
$ws.Range("R448").Value = "Hortaliza"

# Fill in the new row 449 with the new weekly record.
$ws.Range("A449").Value = 11
$ws.Range("B449").Value = "Vega Monumental Concepción"
$ws.Range("C449").Value = "Bíobío"
$ws.Range("D449").Value = 45223
$ws.Range("E449").Value = 8
$ws.Range("F449").Value = 100114013
$ws.Range("G449").Value = "Zanahoria"
$ws.Range("H449").Value = "Sin especificar"
$ws.Range("I449").Value = "Primera"
$ws.Range("J449").Value = 300
$ws.Range("K449").Value = 4000
$ws.Range("L449").Value = 4000
$ws.Range("M449").Value = 4000
$ws.Range("N449").Value = "$/saco 20 kilos"
$ws.Range("O449").Value = "Región de Ñuble"
$ws.Range("P449").Value = 200
$ws.Range("Q449").Value = 20
$ws.Range("R449").Value = "Hortaliza"
